$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '66.240.93'
$ws.Range('E2').Value = '  -0.03%  '

# Row 3
$ws.Range('D3').Value = '3.544.69'
$ws.Range('E3').Value = '  -0.14%  '

# Row 4
$ws.Range('E4').Value = '  -0.10%  '

# Row 5
$ws.Range('D5').Value = '602.76'
$ws.Range('E5').Value = '  -0.38%  '

# Row 6
$ws.Range('D6').Value = '145.85'
$ws.Range('E6').Value = '  +1.24%  '

# Row 7
$ws.Range('D7').Value = '3.545.55'
$ws.Range('E7').Value = '  -0.12%  '

# Row 8
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.26%  '

# Row 9
$ws.Range('D9').Value = '0.496'
$ws.Range('E9').Value = '  +1.29%  '

# Row 10
$ws.Range('E10').Value = '  -1.87%  '

# Row 11
$ws.Range('D11').Value = '7.77'
$ws.Range('E11').Value = '  -0.65%  '

# Row 12
$ws.Range('D12').Value = '0.408'
$ws.Range('E12').Value = '  -1.04%  '

# Row 13
$ws.Range('D13').Value = '4.142.02'
$ws.Range('E13').Value = '  -0.27%  '

# Row 14
$ws.Range('E14').Value = '  -2.51%  '

# Row 15
$ws.Range('D15').Value = '29.02'
$ws.Range('E15').Value = '  -3.37%  '

# Row 16
$ws.Range('D16').Value = '3.550.45'
$ws.Range('E16').Value = '  -0.27%  '

# Row 17
$ws.Range('E17').Value = '  +1.49%  '

# Row 18
$ws.Range('D18').Value = '66.168.55'
$ws.Range('E18').Value = '  -0.29%  '

# Row 19
$ws.Range('E19').Value = '  -3.90%  '

# Row 20
$ws.Range('D20').Value = '6.24'
$ws.Range('E20').Value = '  +1.00%  '

# Row 21
$ws.Range('D21').Value = '14.65'
$ws.Range('E21').Value = '  -0.94%  '

# Row 22
$ws.Range('D22').Value = '417.50'
$ws.Range('E22').Value = '  -3.04%  '

# Row 23
$ws.Range('D23').Value = '0.600'
$ws.Range('E23').Value = '  -1.55%  '

# Row 24
$ws.Range('D24').Value = '77.63'
$ws.Range('E24').Value = '  -2.44%  '

# Row 25
$ws.Range('D25').Value = '3.681.02'
$ws.Range('E25').Value = '  -0.36%  '

# Row 26
$ws.Range('E26').Value = '  +0.16%  '

# Row 27
$ws.Range('E27').Value = '  -2.45%  '

# Row 28
$ws.Range('D28').Value = '9.10'
$ws.Range('E28').Value = '  +0.01%  '

# Row 29
$ws.Range('E29').Value = '  -1.40%  '

# Row 30
$ws.Range('D30').Value = '7.78'
$ws.Range('E30').Value = '  -1.84%  '

# Row 31
$ws.Range('E31').Value = '  -0.02%  '

# Row 32
$ws.Range('D32').Value = '3.539.27'
$ws.Range('E32').Value = '  -0.16%  '

# Row 33
$ws.Range('E33').Value = '  +2.28%  '

# Row 34
$ws.Range('D34').Value = '24.42'
$ws.Range('E34').Value = '  -3.86%  '

# Row 35
$ws.Range('E35').Value = '  +0.01%  '

# Row 36
$ws.Range('D36').Value = '7.56'
$ws.Range('E36').Value = '  -3.27%  '

# Row 37
$ws.Range('E37').Value = '  -10.23%  '

# Row 38
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').Value = '173.98'
$ws.Range('E38').Value = '  -1.14%  '

# Row 39
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').Value = '1.60'
$ws.Range('E39').Value = '  -6.87%  '

# Row 40
$ws.Range('E40').Value = '  -5.65%  '

# Row 41
$ws.Range('D41').Value = '0.0821'
$ws.Range('E41').Value = '  -2.95%  '

# Row 42
$ws.Range('D42').Value = '5.08'
$ws.Range('E42').Value = '  -2.02%  '

# Row 43
$ws.Range('D43').Value = '0.860'
$ws.Range('E43').Value = '  -2.98%  '

# Row 44
$ws.Range('D44').Value = '45.65'
$ws.Range('E44').Value = '  -0.64%  '

# Row 45
$ws.Range('E45').Value = '  -6.09%  '

# Row 46
$ws.Range('D46').Value = '1.00'
$ws.Range('E46').Value = '  -0.02%  '

# Row 47
$ws.Range('D47').Value = '2.40'
$ws.Range('E47').Value = '  -3.91%  '

# Row 48
$ws.Range('D48').Value = '7.12'
$ws.Range('E48').Value = '  -0.01%  '

# Row 49
$ws.Range('D49').Value = '22.76'
$ws.Range('E49').Value = '  -1.89%  '

# Row 50
$ws.Range('E50').Value = '  -7.74%  '

# Row 51
$ws.Range('D51').Value = '23.17'
$ws.Range('E51').Value = '  -7.75%  '
